# Sprint 2 Daily Scrum - add the entry for day 2021-02-03 (row block 25-29),
# mirroring the layout of the previous day blocks (header row + one row per
# team member, in the workbook's usual order: Anibal, Julio, Pedro, Vitor).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (25), cloned from the previous header row (19) --------------
$ws.Range("A19:D19").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122)
$ws.Rows.Item(25).RowHeight = 163.5

$ws.Range("A25").Value = 44230
$ws.Range("B25").Value = "1- Whats`n been done"
$ws.Range("C25").Value = "2- what `nwill we do"
$ws.Range("D25").Value = "3- difficulties`n encountered"

# --- Anibal row (26), cloned from row 20 -------------------------------------
$ws.Range("A20:D20").Copy()
$ws.Range("A26:D26").PasteSpecial(-4122)
$ws.Rows.Item(26).RowHeight = 163.5

$ws.Range("A26").Value = "Anibal"
$ws.Range("B26").Value = "implementar classes registarColaboradorcontroller e areaGestorUI"
$ws.Range("C26").Value = "terminar a areaColaboradorUI"
$ws.Range("D26").Value = "NADA A APONTAR"

# --- Julio row (27), cloned from row 21 --------------------------------------
$ws.Range("A21:D21").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122)
$ws.Rows.Item(27).RowHeight = 163.5

$ws.Range("A27").Value = "Julio"
$ws.Range("B27").Value = "INICIADO INTERFACES GRAFICAS UC2E E IMPLEMENTACAO UC6"
$ws.Range("C27").Value = "IMPLEMENTAR UI UC6"
$ws.Range("D27").Value = "NADA A APONTAR"

# --- Pedro row (28), cloned from row 22 --------------------------------------
$ws.Range("A22:D22").Copy()
$ws.Range("A28:D28").PasteSpecial(-4122)
$ws.Rows.Item(28).RowHeight = 163.5

$ws.Range("A28").Value = "Pedro"
$ws.Range("B28").Value = "TERMINAR LIGACAO COM A API`nTERMINAR SERIALIZACAO DADOS"
$ws.Range("C28").Value = "plano a - revisao javafx`ncaso haja tempo trabalhar nas classes UI`nterminar uc1"
$ws.Range("D28").Value = "NADA A APONTAR"

# --- Vitor row (29), cloned from row 23 --------------------------------------
$ws.Range("A23:D23").Copy()
$ws.Range("A29:D29").PasteSpecial(-4122)
$ws.Rows.Item(29).RowHeight = 163.5

$ws.Range("A29").Value = "Vitor"
$ws.Range("B29").Value = "CRIACAO CATEGORIATAREFA E CARACTERIZACAOCOMPTEC E RESPETIVO CONTROLLER`nPESQUISA CSS"
$ws.Range("C29").Value = "FINALIZAR TDD CLASSES CRIADAS"
$ws.Range("D29").Value = "NADA A APONTAR"

# --- Selection / view state, matching the author's final cursor position ----
$ws.Range("D29").Select()
